$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.005.10'
$ws.Range('E2').Value = '  -2.33%  '
$ws.Range('D3').Value = '3.118.27'
$ws.Range('E3').Value = '  -0.90%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '594.04'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -2.39%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '136.05'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -5.63%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '3.112.44'
$ws.Range('E8').Value = '  -0.97%  '
$ws.Range('E9').Value = '  -2.23%  '
$ws.Range('E10').Value = '  -4.02%  '
$ws.Range('E11').Value = '  -3.05%  '
$ws.Range('E12').Value = '  -2.84%  '
$ws.Range('E13').Value = '  -4.87%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '34.09'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -3.76%  '
$ws.Range('D15').Value = '3.626.51'
$ws.Range('E15').Value = '  -0.98%  '
$ws.Range('E16').Value = '  +2.56%  '
$ws.Range('D17').Value = '62.950.50'
$ws.Range('E17').Value = '  -2.19%  '
$ws.Range('D18').Value = '3.117.77'
$ws.Range('E18').Value = '  -0.93%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.66'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -2.86%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '472.37'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.83%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.08'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -4.30%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.696'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -2.71%  '
$ws.Range('E23').Value = '  -0.52%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '85.97'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.52%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '12.91'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -3.57%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('E27').Value = '  -1.92%  '
$ws.Range('B28').Value = 'NEARProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.93'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -3.62%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.88'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -6.80%  '
$ws.Range('E30').Value = '  +1.52%  '
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '26.66'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.53%  '
$ws.Range('E33').Value = '  -6.63%  '
$ws.Range('E34').Value = '  -4.83%  '
$ws.Range('E35').Value = '  -3.15%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.79'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -3.17%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '52.00'
$ws.Range('D37').Style = "Normal"
$ws.Range('D38').Value = '0.0₃0697'
$ws.Range('E38').Value = '  -9.35%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0386'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.84%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '417.12'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -6.93%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '8.18'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.78%  '
$ws.Range('D42').Value = '2.900.15'
$ws.Range('E42').Value = '  +0.89%  '
$ws.Range('E43').Value = '  -11.56%  '
$ws.Range('E44').Value = '  -5.99%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.266'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.52%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.11'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -5.63%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '25.46'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.89%  '
$ws.Range('E49').Value = '  -0.65%  '
$ws.Range('E50').Value = '  -6.81%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '119.65'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.04%  '
